# Write time reformatting function
#
# Populate the (empty) document body with a list of times, reformatted
# from mm:ss.hh style into total seconds (one value per paragraph),
# inserted just before the final section break.

$d = $word.ActiveDocument

$times = @(
    "0.69", "2.50", "9.71", "12.09", "16.29",
    "20.15", "24.27", "26.93", "32.28", "37.94",
    "45.08", "49.10", "53.29", "60.18", "70.19",
    "76.01", "83.43", "89.33", "97.28", "111.97",
    "115.91", "123.44", "125.72", "130.98", "143.84",
    "147.35", "151.34", "153.06", "181.01", "184.06",
    "191.89", "194.72", "202.78", "210.62", "215.22",
    "222.63", "230.81", "235.47", "241.35", "250.08",
    "262.59", "266.00", "276.28", "282.34", "283.67"
)

for ($i = 0; $i -lt $times.Count; $i++) {
    $d.Content.InsertAfter($times[$i])
    if ($i -lt $times.Count - 1) {
        $d.Content.InsertParagraphAfter()
    }
}
